$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.314.29'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '1.888.97'
$ws.Range("E3").Value = '  -1.18%  '
$ws.Range("E4").Value = '  -0.91%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.99'
$ws.Range("E5").Value = '  -2.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.690'
$ws.Range("E6").Value = '  -4.77%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.98%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '43.99'
$ws.Range("E8").Value = '  +8.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.351'
$ws.Range("E9").Value = '  -3.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '51.77'
$ws.Range("E10").Value = '  -0.98%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0737'
$ws.Range("E11").Value = '  -3.91%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0971'
$ws.Range("E12").Value = '  -1.75%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '13.14'
$ws.Range("E13").Value = '  +3.25%  '
$ws.Range("D14").Value = '2.161.77'
$ws.Range("E14").Value = '  -1.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.724'
$ws.Range("E15").Value = '  +0.67%  '
$ws.Range("D16").Value = '1.909.39'
$ws.Range("E16").Value = '  +0.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.92'
$ws.Range("E17").Value = '  -0.15%  '
$ws.Range("D18").Value = '35.260.29'
$ws.Range("E18").Value = '  +0.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.11'
$ws.Range("E19").Value = '  -1.76%  '
$ws.Range("D20").Value = '0.0₃0821'
$ws.Range("E20").Value = '  -3.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '245.85'
$ws.Range("E21").Value = '  +0.86%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.83'
$ws.Range("E22").Value = '  -1.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.97'
$ws.Range("E23").Value = '  -2.27%  '
$ws.Range("E24").Value = '  -0.95%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.53'
$ws.Range("E25").Value = '  +6.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.19'
$ws.Range("E26").Value = '  -11.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '165.34'
$ws.Range("E27").Value = '  -0.81%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.50'
$ws.Range("E28").Value = '  -2.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.33'
$ws.Range("E30").Value = '  -4.04%  '
$ws.Range("D31").Value = '4.128.49'
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.78'
$ws.Range("E32").Value = '  +8.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.27'
$ws.Range("E33").Value = '  -1.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0580'
$ws.Range("E34").Value = '  -0.46%  '
$ws.Range("E35").Value = '  +0.43%  '
$ws.Range("E36").Value = '  -0.95%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.855'
$ws.Range("E37").Value = '  -6.93%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.00'
$ws.Range("E38").Value = '  -1.88%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.57'
$ws.Range("E39").Value = '  -21.62%  '
$ws.Range("E40").Value = '  +0.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '97.80'
$ws.Range("E41").Value = '  +1.14%  '
$ws.Range("E42").Value = '  +3.23%  '
$ws.Range("E43").Value = '  -2.25%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.10'
$ws.Range("E44").Value = '  -2.20%  '
$ws.Range("D45").Value = '1.291.19'
$ws.Range("E45").Value = '  -3.58%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0806'
$ws.Range("E47").Value = '  +8.10%  '
$ws.Range("E48").Value = '  -0.86%  '
$ws.Range("E49").Value = '  -1.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '12.09'
$ws.Range("E50").Value = '  +0.38%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.42'
$ws.Range("E51").Value = '  -5.17%  '
